$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M138").ClearContents()
$ws.Range("H138").Value = 3965
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3965
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 11895
$ws.Range("N138").Value = -22175

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9260473
$ws.Range("I61").Value = 12501202
$ws.Range("J61").Value = 1247.0714
$ws.Range("K61").Value = 12501202
$ws.Range("L61").Value = 1247.0714
$ws.Range("M61").Value = -12500990
$ws.Range("N61").Value = -1671.0714
$ws.Range("H132").Value = 5103773.5
$ws.Range("I132").Value = 7577318
$ws.Range("J132").Value = 2087.5625
$ws.Range("K132").Value = 22731954
$ws.Range("L132").Value = 6262.6875
$ws.Range("M132").Value = -22729424
$ws.Range("N132").Value = -11322.6875
$ws.Range("H136").Value = 9260473
$ws.Range("I136").Value = 12501202
$ws.Range("J136").Value = 1247.0714
$ws.Range("K136").Value = 37503606
$ws.Range("L136").Value = 3741.2142
$ws.Range("M136").Value = -37501056
$ws.Range("N136").Value = -8841.2142

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1866.1637
$ws.Range("I134").Value = 1137.925
$ws.Range("J134").Value = 3808.1333
$ws.Range("K134").Value = 3413.775
$ws.Range("L134").Value = 11424.3999
$ws.Range("M134").Value = -878.7749999999996
$ws.Range("N134").Value = -16494.3999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7580149.5
$ws.Range("I31").Value = 5124.9355
$ws.Range("K31").Value = 5124.9355
$ws.Range("M31").Value = -4829.9355
$ws.Range("H34").Value = 7580149.5
$ws.Range("I34").Value = 5124.9355
$ws.Range("K34").Value = 5124.9355
$ws.Range("M34").Value = -4922.9355
$ws.Range("H94").Value = 2935
$ws.Range("I94").Value = 2073.8333
$ws.Range("J94").Value = 3451.7
$ws.Range("K94").Value = 2073.8333
$ws.Range("L94").Value = 3451.7
$ws.Range("M94").Value = -1622.8333
$ws.Range("N94").Value = -4353.7
$ws.Range("H129:L129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130:L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131:L131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:L133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H137:L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138:L138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141:L141").ClearContents()
$ws.Range("N141").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 13626.818
$ws.Range("I120").Value = 7030
$ws.Range("J120").Value = 14286.5
$ws.Range("K120").Value = 21090
$ws.Range("L120").Value = 42859.5
$ws.Range("M120").Value = -16252
$ws.Range("N120").Value = -52535.5
$ws.Range("H121").Value = 900.7
$ws.Range("I121").Value = 266
$ws.Range("J121").Value = 1112.2667
$ws.Range("K121").Value = 798
$ws.Range("L121").Value = 3336.800099999999
$ws.Range("M121").Value = 512
$ws.Range("N121").Value = -5956.800099999999
$ws.Range("H122").Value = 1844.375
$ws.Range("I122").Value = 1959.3334
$ws.Range("J122").Value = 1499.5
$ws.Range("K122").Value = 17634.0006
$ws.Range("L122").Value = 13495.5
$ws.Range("M122").Value = -15184.0006
$ws.Range("N122").Value = -18395.5
$ws.Range("H123").Value = 6416.6665
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 6416.6665
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 19249.9995
$ws.Range("N123").Value = -24149.9995
$ws.Range("H124").Value = 7703.3335
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 7703.3335
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 23110.0005
$ws.Range("N124").Value = -32930.00049999999
$ws.Range("H125").Value = 3148
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 3640
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 10920
$ws.Range("M125").Value = -1080
$ws.Range("N125").Value = -20760
$ws.Range("H126").Value = 5391.6665
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5391.6665
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 16174.9995
$ws.Range("N126").Value = -26054.9995
$ws.Range("H127").Value = 818.6
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 818.6
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 2455.8
$ws.Range("N127").Value = -12375.8
$ws.Range("H128").Value = 300000
$ws.Range("I128").Value = 300000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 900000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -895020
$ws.Range("H129").Value = 2989.68
$ws.Range("I129").Value = 2916.6667
$ws.Range("J129").Value = 2999.6365
$ws.Range("K129").Value = 8750.000100000001
$ws.Range("L129").Value = 8998.9095
$ws.Range("M129").Value = -3750.000100000001
$ws.Range("N129").Value = -18998.9095
$ws.Range("H130").Value = 5497.143
$ws.Range("I130").Value = 2353.3333
$ws.Range("J130").Value = 6354.5454
$ws.Range("K130").Value = 7059.999899999999
$ws.Range("L130").Value = 19063.6362
$ws.Range("M130").Value = -2039.999899999999
$ws.Range("N130").Value = -29103.6362
$ws.Range("H131").Value = 858.38
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 858.38
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2575.14
$ws.Range("N131").Value = -12655.14
$ws.Range("H132").Value = 1171.2
$ws.Range("I132").Value = 838.375
$ws.Range("J132").Value = 2502.5
$ws.Range("K132").Value = 7545.375
$ws.Range("L132").Value = 22522.5
$ws.Range("M132").Value = -5015.375
$ws.Range("N132").Value = -27582.5
$ws.Range("H133").Value = 4992
$ws.Range("I133").Value = 2235
$ws.Range("J133").Value = 8142.857
$ws.Range("K133").Value = 6705
$ws.Range("L133").Value = 24428.571
$ws.Range("M133").Value = -1645
$ws.Range("N133").Value = -34548.571
$ws.Range("H134").Value = 3602.8948
$ws.Range("I134").Value = 1783.0435
$ws.Range("J134").Value = 6393.3335
$ws.Range("K134").Value = 5349.1305
$ws.Range("L134").Value = 19180.0005
$ws.Range("M134").Value = -279.1305000000002
$ws.Range("N134").Value = -29320.0005
$ws.Range("H136").Value = 2575.4546
$ws.Range("I136").Value = 1203.3334
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 3610.0002
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = 1489.9998
$ws.Range("N136").Value = -36450
$ws.Range("H137").Value = 7292.1665
$ws.Range("I137").Value = 5986
$ws.Range("J137").Value = 7794.5386
$ws.Range("K137").Value = 17958
$ws.Range("L137").Value = 23383.6158
$ws.Range("M137").Value = -12858
$ws.Range("N137").Value = -33583.6158
$ws.Range("H138").Value = 10298
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 10298
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 30894
$ws.Range("N138").Value = -41174
$ws.Range("H139").Value = 2205.9285
$ws.Range("I139").Value = 1570.8334
$ws.Range("J139").Value = 6016.5
$ws.Range("K139").Value = 4712.5002
$ws.Range("L139").Value = 18049.5
$ws.Range("M139").Value = 427.4997999999996
$ws.Range("N139").Value = -28329.5
$ws.Range("H140").Value = 3411.6843
$ws.Range("I140").Value = 1614.5385
$ws.Range("J140").Value = 7305.5
$ws.Range("K140").Value = 4843.6155
$ws.Range("L140").Value = 21916.5
$ws.Range("M140").Value = 336.3845000000001
$ws.Range("N140").Value = -32276.5
$ws.Range("H141").Value = 5750.5264
$ws.Range("I141").Value = 4074.1177
$ws.Range("J141").Value = 20000
$ws.Range("K141").Value = 12222.3531
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = -7042.3531
$ws.Range("N141").Value = -70360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3819.65
$ws.Range("I126").Value = 2517.3333
$ws.Range("J126").Value = 5259.0527
$ws.Range("K126").Value = 7551.999899999999
$ws.Range("L126").Value = 15777.1581
$ws.Range("M126").Value = -5081.999899999999
$ws.Range("N126").Value = -20717.1581
$ws.Range("H132").Value = 2546.9553
$ws.Range("I132").Value = 1822.3529
$ws.Range("J132").Value = 4856.625
$ws.Range("K132").Value = 5467.0587
$ws.Range("L132").Value = 14569.875
$ws.Range("M132").Value = -2937.0587
$ws.Range("N132").Value = -19629.875
$ws.Range("H138").Value = 59950
$ws.Range("J138").Value = 59950
$ws.Range("L138").Value = 59950
$ws.Range("N138").Value = -70230

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1297.1666
$ws.Range("I93").Value = 1141.5714
$ws.Range("J93").Value = 1515
$ws.Range("K93").Value = 1141.5714
$ws.Range("L93").Value = 1515
$ws.Range("M93").Value = 106.4286
$ws.Range("N93").Value = -4011
$ws.Range("H122").Value = 4234.516
$ws.Range("I122").Value = 4346.1904
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 13038.5712
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -10588.5712
$ws.Range("N122").Value = -16900
$ws.Range("H136").Value = 11370534
$ws.Range("I136").Value = 16668883
$ws.Range("J136").Value = 16929.285
$ws.Range("K136").Value = 50006649
$ws.Range("L136").Value = 50787.855
$ws.Range("M136").Value = -50004099
$ws.Range("N136").Value = -55887.855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1479.0793
$ws.Range("I132").Value = 1217.2885
$ws.Range("K132").Value = 3651.8655
$ws.Range("M132").Value = -1121.8655
$ws.Range("H136").Value = 585.4912
$ws.Range("I136").Value = 539.46
$ws.Range("J136").Value = 914.2857
$ws.Range("K136").Value = 1618.38
$ws.Range("L136").Value = 2742.8571
$ws.Range("M136").Value = 931.6199999999999
$ws.Range("N136").Value = -7842.8571
